$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EICHER")

# Row 7 (summary row at top: F=close/LTP ref row, G=high, H=low, I=ltp, J=prev)
$ws.Range("F7").Value = 3855
$ws.Range("G7").Value = 4007.9
$ws.Range("H7").Value = 3847.55
$ws.Range("I7").Value = 3957.95
$ws.Range("J7").Value = 3912.4

# Row 9
$ws.Range("G9").Value = 3908.4
$ws.Range("H9").Value = 3847.55
$ws.Range("I9").Value = 3864.3

# Row 10
$ws.Range("G10").Value = 3945.25
$ws.Range("H10").Value = 3862.4
$ws.Range("I10").Value = 3922

# Row 11
$ws.Range("G11").Value = 3973.4
$ws.Range("H11").Value = 3919.95
$ws.Range("I11").Value = 3942.75

# Row 12
$ws.Range("G12").Value = 3967.9
$ws.Range("H12").Value = 3925
$ws.Range("I12").Value = 3951.95

# Row 13
$ws.Range("G13").Value = 3951.95
$ws.Range("H13").Value = 3932.1
$ws.Range("I13").Value = 3944.35

# Row 14
$ws.Range("G14").Value = 3961.15
$ws.Range("H14").Value = 3942.6
$ws.Range("I14").Value = 3957.65

# Row 15
$ws.Range("G15").Value = 3967.1
$ws.Range("H15").Value = 3945
$ws.Range("I15").Value = 3952.8

# Row 16
$ws.Range("G16").Value = 3965.95
$ws.Range("H16").Value = 3946.15
$ws.Range("I16").Value = 3964.9

# Row 17
$ws.Range("G17").Value = 3970
$ws.Range("H17").Value = 3955.6
$ws.Range("I17").Value = 3969

# Row 18
$ws.Range("G18").Value = 4007.9
$ws.Range("H18").Value = 3968
$ws.Range("I18").Value = 3997

# Row 19
$ws.Range("G19").Value = 3997.05
$ws.Range("H19").Value = 3975.05
$ws.Range("I19").Value = 3988

# Row 20
$ws.Range("G20").Value = 3990
$ws.Range("H20").Value = 3958.8
$ws.Range("I20").Value = 3962.4

# Row 21
$ws.Range("G21").Value = 3978.4
$ws.Range("H21").Value = 3949.45
$ws.Range("I21").Value = 3963.85
